$p = $ppt.ActivePresentation

# The existing slide 1 ("Projet Booker") is a Title Slide (ctrTitle + subTitle
# placeholders). The new closing "FIN" slide uses the very same layout, so the
# simplest and most faithful way to build it is to duplicate slide 1 and then
# change its contents, which preserves the placeholder/style structure exactly.
$src = $p.Slides.Item(1)
$newSlide = $src.Duplicate()

# Move the freshly duplicated slide (currently right after slide 1) to the end
# of the deck, as the 13th and last slide.
$newSlide.MoveTo($p.Slides.Count)

# Title placeholder: "Projet Booker" -> "FIN"
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "FIN"

# Subtitle placeholder: shrink text to fit (-> <a:normAutofit/>) and make it
# three empty paragraphs, like the authored slide.
$subTitleShape = $newSlide.Shapes.Item(2)
$subTitleShape.TextFrame.AutoSize = 2

$subTitleRange = $subTitleShape.TextFrame.TextRange
$subTitleRange.Text = "x`rx`rx"
$subTitleRange.Characters(1, 1).Text = ""
$subTitleRange.Characters(2, 1).Text = ""
$subTitleRange.Characters(3, 1).Text = ""
